$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K9").Value = 131198.18
$ws.Range("M10").Value = 439284.65
$ws.Range("O10").Value = 175224.62
$ws.Range("O21").Value = 679963.07
$ws.Range("O29").Value = 202098
$ws.Range("O30").Value = 18844.7
